$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for every data row (2-395)
#    from 2023-09-13 (45182) to 2023-09-15 (45184).
for ($r = 2; $r -le 395; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# 2. Row 395 picks up an explicit custom row height (matches the other rows).
$ws.Rows.Item(395).RowHeight = 15

# 3. Append the new record as row 396.
$ws.Cells.Item(396, 1).Value = "A 43264-2023"

$ws.Cells.Item(396, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(396, 2).Value = 45183

$ws.Cells.Item(396, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(396, 3).Value = 45184

$ws.Cells.Item(396, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(396, 5).Value = "HALMSTAD"

$ws.Cells.Item(396, 7).Value = 2.3
$ws.Cells.Item(396, 8).Value = 0
$ws.Cells.Item(396, 9).Value = 0
$ws.Cells.Item(396, 10).Value = 0
$ws.Cells.Item(396, 11).Value = 0
$ws.Cells.Item(396, 12).Value = 0
$ws.Cells.Item(396, 13).Value = 0
$ws.Cells.Item(396, 14).Value = 0
$ws.Cells.Item(396, 15).Value = 0
$ws.Cells.Item(396, 16).Value = 0
$ws.Cells.Item(396, 17).Value = 0

# Column R keeps the same wrap-text styling used throughout the sheet,
# but stays empty (no species names recorded for this entry).
$ws.Range("R396").WrapText = $true
